# Daily attendance processing - 2025-11-24 20:29:35
#
# Normalises the "Recorded By" column (G) on the session-analysis sheet:
# each cell holds a comma-separated list of recorders (e.g. "System,
# dnasr281@gmail.com"); re-sort that list alphabetically (case-insensitive)
# so entries like "System" no longer always lead.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # column G = "Recorded By"
    $text = $cell.Text

    if ([string]::IsNullOrEmpty($text)) { continue }

    $parts = $text -split ", "
    if ($parts.Count -le 1) { continue }

    $sorted = ($parts | Sort-Object) -join ", "

    if ($sorted -ne $text) {
        $cell.Value = $sorted
    }
}
